# Apply "2 down, 2 to go" parameter-table update:
#   - row 10 (name column) is renamed from "scale_rh" to "f" on all three
#     sheets (the "Scaling factor for fluxes" description / units stay put)
#   - a new row 11 is appended on all three sheets for the "gR" parameter
#     ("Base autotrophic respiration rate", units "g C m^-2 day-1")
#   - on the "ranges" sheet the scale_rh/f row's min/max values are updated
#     (0 -> 1, 100000000 -> 10) and the new gR row gets value=1, min=0, max=10
#   - the leftover manual formatting on the "incubation" sheet's row 10 is
#     cleared so it matches the other rows again

$wb = $excel.ActiveWorkbook

$incubation = $wb.Worksheets.Item("incubation")
$field      = $wb.Worksheets.Item("field")
$ranges     = $wb.Worksheets.Item("ranges")

# --- rename scale_rh -> f on every sheet (description/units unchanged) ---
$incubation.Range("A10").Value = "f"
$field.Range("A10").Value = "f"
$ranges.Range("A10").Value = "f"

# --- incubation: drop the stray cell formatting that had crept onto A10:C10 ---
$incubation.Range("A10:C10").ClearFormats()

# --- new gR row on "incubation" (bool flag columns D:H) ---
$incubation.Range("C11").Value = "g C m^-2 day-1"
$incubation.Range("A11").Value = "gR"
$incubation.Range("B11").Value = "Base autotrophic respiration rate"
$incubation.Range("D11").Value = $false
$incubation.Range("E11").Value = $false
$incubation.Range("F11").Value = $false
$incubation.Range("G11").Value = $false
$incubation.Range("H11").Value = $false

# --- new gR row on "field" (bool flag columns D:H) ---
$field.Range("C11").Value = "g C m^-2 day-1"
$field.Range("A11").Value = "gR"
$field.Range("B11").Value = "Base autotrophic respiration rate"
$field.Range("D11").Value = $false
$field.Range("E11").Value = $false
$field.Range("F11").Value = $false
$field.Range("G11").Value = $false
$field.Range("H11").Value = $false

# --- ranges sheet: update f's min/max, then add the gR row (value/min/max) ---
$ranges.Range("E10").Value = 1
$ranges.Range("F10").Value = 10

$ranges.Range("C11").Value = "g C m^-2 day-1"
$ranges.Range("A11").Value = "gR"
$ranges.Range("B11").Value = "Base autotrophic respiration rate"
$ranges.Range("D11").Value = 1
$ranges.Range("E11").Value = 0
$ranges.Range("F11").Value = 10

# --- selections: the newly-added row is now what's selected on each sheet ---
$field.Range("A11:B11").Select()

$ranges.Range("A11:B11").Select()

$incubation.Range("A11:B11").Select()
$incubation.Activate()
